# Update the "Integer min" value for rule R20 (row 10, column C) from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
